# Generate Report for Archive
#
# The CI-generated localization-status report is refreshed: the rows for
# 8905827f-75f3-4725-98c1-1e1795ff5034 and 2846089b-1542-4c33-a715-291019ee01ea
# swap their relative order (2846089b now appears before 8905827f in rows 4/5
# of every sheet), their "Latest *" timestamps/files follow the swap, and the
# status of both rows becomes "In Translation".

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink {
    param($ws, $cellAddr, $value)
    $ws.Range($cellAddr).Value = $value
    $target = '$' + $cellAddr.Substring(0,1) + '$' + $cellAddr.Substring(1)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $target) {
            $h.TextToDisplay = $value
        }
    }
}

# ----------------------------------------------------------------------
# Sheet "Overview": A=File Name, B=Path And Name (hyperlink), C=Extension,
# D=Publish URL, E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $wsOverview "A4" "2846089b-1542-4c33-a715-291019ee01ea.md"
Set-CellAndHyperlink $wsOverview "B4" "e2e\2846089b-1542-4c33-a715-291019ee01ea.md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-08-26 16:42:09"

Set-CellAndHyperlink $wsOverview "A5" "8905827f-75f3-4725-98c1-1e1795ff5034.md"
Set-CellAndHyperlink $wsOverview "B5" "e2e\8905827f-75f3-4725-98c1-1e1795ff5034.md"
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("G5").Value = "2016-08-26 16:41:30"

# ----------------------------------------------------------------------
# Sheet "zh-cn": A=Source File Name (hyperlink), C=Status,
# G=Latest Handoff File, H=Latest Handoff Datetime
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $wsZhCn "A4" "2846089b-1542-4c33-a715-291019ee01ea.md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("G4").Value = "2846089b-1542-4c33-a715-291019ee01ea.7d99c404bec278dea52c83eaada485a1f65209c9.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-26 16:41:57"

Set-CellAndHyperlink $wsZhCn "A5" "8905827f-75f3-4725-98c1-1e1795ff5034.md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("G5").Value = "8905827f-75f3-4725-98c1-1e1795ff5034.91c313994e076b9d29fda557407158aae7e579f0.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-26 16:41:25"

# ----------------------------------------------------------------------
# Sheet "de-de": A=Source File Name (hyperlink), C=Status,
# G=Latest Handoff File, H=Latest Handoff Datetime
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $wsDeDe "A4" "2846089b-1542-4c33-a715-291019ee01ea.md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("G4").Value = "2846089b-1542-4c33-a715-291019ee01ea.7d99c404bec278dea52c83eaada485a1f65209c9.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-26 16:42:09"

Set-CellAndHyperlink $wsDeDe "A5" "8905827f-75f3-4725-98c1-1e1795ff5034.md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("G5").Value = "8905827f-75f3-4725-98c1-1e1795ff5034.91c313994e076b9d29fda557407158aae7e579f0.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-26 16:41:30"
